# Time log + project profile update
# Fills in the first week of logged time entries (dates, durations, activity
# notes) and adds the "Total Time:" summary formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Date 2/4/2020, 1 hour, "logo" ---
$ws.Range("A4").Value = 43865
$ws.Range("A4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "logo"

# --- Row 5: Date 2/5/2020, 0.25 hour, "added logo " ---
$ws.Range("A5").Value = 43866
$ws.Range("A5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = "added logo "

# --- Row 6: Date 2/6/2020, 0.25 hour, "updated logo sizes" ---
$ws.Range("A6").Value = 43867
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("D6").Value = 0.25
$ws.Range("E6").Value = "updated logo sizes"

# --- Row 7: Date 2/10/2020, 2 hours, "about me page" ---
$ws.Range("A7").Value = 43871
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "about me page"

# --- Row 9: Date 2/12/2020, 0.25 hour, "Fixed experience text" ---
# (entered before row 8, matching the original authoring order)
$ws.Range("A9").Value = 43873
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("D9").Value = 0.25
$ws.Range("E9").Value = "Fixed experience text"

# --- Row 8: Date 2/11/2020, Start 7:00 PM, Stop 8:00 PM, "Built connect 4 board" ---
$ws.Range("A8").Value = 43872
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = 0.79166666666666663
$ws.Range("B8").NumberFormat = "h:mm AM/PM"
$ws.Range("C8").Value = 0.83333333333333337
$ws.Range("C8").NumberFormat = "h:mm AM/PM"
$ws.Range("D8").Formula = '=IF(OR(ISBLANK(B8),ISBLANK(C8)),"",(C8-B8)*24)'
$ws.Range("E8").Value = "Built connect 4 board"

# --- Fill the rest of the Delta column (D10:D31) with the cleaned-up,
#     #REF!-free formula, as one shared-formula fill-down. ---
$ws.Range("D10:D31").Formula = '=IF(OR(ISBLANK(B10),ISBLANK(C10)),"",(C10-B10)*24)'

# Give row 10's date/start-time cells their number formats too (no values).
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("B10").NumberFormat = "h:mm"

# --- H68:H92 "Invalid stop time" checks collapse into one shared formula ---
$ws.Range("H68:H92").Formula = '=IF(D68<0,"<-- Invalid stop time","")'

# --- Row 2: Total Time label + sum formula ---
$ws.Range("C2").Value = "Total Time:"
$ws.Range("D2").Formula = "=SUM(D4:D44)"

# --- Selection moves to A10:B10 (next unfilled entry) ---
$ws.Range("A10:B10").Select()
